# Update the DQ_Metrics sheet:
#  - rename / re-derive several metric columns (C..S)
#  - drop the now-unused trailing "_no_py" raw-count columns (old T:AA)
#  - refresh the row-2 data values to match the new column semantics

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DQ_Metrics")

# Drop the eight trailing helper columns (old T:AA) - this also shrinks
# the sheet dimension from A1:AA2 down to A1:S2.
$ws.Range("T1:AA2").EntireColumn.Delete()

# New header row (columns C..S)
$headers = @{
    "C1" = "missing_item_rate"
    "D1" = "missing_value_rate"
    "E1" = "orphaCoding_completeness_rate"
    "F1" = "outlier_rate"
    "G1" = "orphaCoding_plausibility_rate"
    "H1" = "rdCase_unambiguity_rate"
    "I1" = "duplication_rate"
    "J1" = "tracerCase_rel_py_ipat"
    "K1" = "unambigous_rdCase_rel_py_ipat"
    "L1" = "orphaCase_rel_py_ipat"
    "M1" = "orphaCoding_no_py"
    "N1" = "orphaCase_no_py"
    "O1" = "unambigous_rdCase_no_py"
    "P1" = "rdCase_no_py"
    "Q1" = "case_no_py"
    "R1" = "patient_no_py"
    "S1" = "case_no_py_ipat"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# Refreshed data row (row 2, columns C..S)
$values = @{
    "C2" = 0
    "D2" = 40.49
    "F2" = 6.25
    "I2" = 3.7
    "J2" = 0.04
    "K2" = 0.07000000000000001
    "L2" = 0.15
    "M2" = 15
    "N2" = 15
    "O2" = 7
    "P2" = 17
    "Q2" = 26
    "R2" = 25
    "S2" = 10000
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
